$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.18370258808136
$ws.Range("B1").Value = 2.313625812530518
$ws.Range("C1").Value = 4.902408599853516
$ws.Range("D1").Value = 2.510704278945923
$ws.Range("E1").Value = 1.219971895217896
